$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Description cell (row 3, column 2): split narrative into new phrasing ---
$descCell = $t.Cell(3, 2)

# Insert "Server for " right after "sent to "
$f1 = $descCell.Range
$f1.Find.Execute("sent to ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p1 = $d.Range($f1.End, $f1.End)
$p1.InsertAfter("Server for ")

# Insert " with Display" right after "shown to the user"
$f2 = $descCell.Range
$f2.Find.Execute("shown to the user", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p2 = $d.Range($f2.End, $f2.End)
$p2.InsertAfter(" with Display")

# Replace "that will be used to get the " with "and "
$f3 = $descCell.Range
$f3.Find.Execute("that will be used to get the ", $true, $false, $false, $false, $false, $true, 1, $false, "and ", 2)

# Replace trailing "related to that user." with "related to that user to store in the Session Storage" (drop period)
$f4 = $descCell.Range
$f4.Find.Execute("related to that user.", $true, $false, $false, $false, $false, $true, 1, $false, "related to that user to store in the Session Storage", 2)

# --- Elements sub-table rows ---
# Row 6 col2: "External Authentication Server - BYUI" -> "1.1.1.1.2 Display"
$c6 = $t.Cell(6, 2)
$r6 = $c6.Range
$r6.Find.Execute("External Authentication Server - BYUI", $true, $false, $false, $false, $false, $true, 1, $false, "1.1.1.1.2 Display", 2)

# Row 7 col2: empty -> "1.1.2 Server"
$c7 = $t.Cell(7, 2)
$r7 = $c7.Range
$r7.InsertBefore("1.1.2 Server")

# Row 8 col2: empty -> "Session Storage"
$c8 = $t.Cell(8, 2)
$r8 = $c8.Range
$r8.InsertBefore("Session Storage")

Write-Host "Edit complete"
